$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Odds refresh (COLOMBIA - PRIMERA A fixtures, rows 2-8 before the deletion below)
# Row 2
$ws.Range("G2").Value2 = 2.9
$ws.Range("H2").Value2 = 3
$ws.Range("I2").Value2 = 2.6
$ws.Range("J2").Value2 = 3.75
$ws.Range("L2").Value2 = 3.4
$ws.Range("M2").Value2 = 1.11
$ws.Range("N2").Value2 = 6.5
$ws.Range("W2").Value2 = 7
$ws.Range("X2").Value2 = 13
$ws.Range("Y2").Value2 = 12
$ws.Range("Z2").Value2 = 29
$ws.Range("AA2").Value2 = 29
$ws.Range("AH2").Value2 = 6.5
$ws.Range("AI2").Value2 = 11
$ws.Range("AJ2").Value2 = 11
$ws.Range("AK2").Value2 = 26
$ws.Range("AL2").Value2 = 26
$ws.Range("AN2").Value2 = 4.75
$ws.Range("AO2").Value2 = 19
$ws.Range("AS2").Value2 = 301
$ws.Range("AW2").Value2 = 4.5
$ws.Range("BB2").Value2 = 301

# Row 3
$ws.Range("G3").Value2 = 2
$ws.Range("H3").Value2 = 3.1
$ws.Range("I3").Value2 = 4.2
$ws.Range("J3").Value2 = 2.88
$ws.Range("L3").Value2 = 4.75
$ws.Range("X3").Value2 = 8.5
$ws.Range("Y3").Value2 = 9.5
$ws.Range("Z3").Value2 = 17
$ws.Range("AA3").Value2 = 19
$ws.Range("AE3").Value2 = 19
$ws.Range("AH3").Value2 = 9
$ws.Range("AI3").Value2 = 19
$ws.Range("AJ3").Value2 = 15
$ws.Range("AL3").Value2 = 41
$ws.Range("AM3").Value2 = 51
$ws.Range("AN3").Value2 = 3.75
$ws.Range("AO3").Value2 = 12
$ws.Range("AP3").Value2 = 26
$ws.Range("AR3").Value2 = 67
$ws.Range("AX3").Value2 = 23
$ws.Range("BB3").Value2 = 351

# Row 4
$ws.Range("G4").Value2 = 4
$ws.Range("I4").Value2 = 2.1

# Row 5
$ws.Range("G5").Value2 = 1.27
$ws.Range("H5").Value2 = 4.75
$ws.Range("I5").Value2 = 15
$ws.Range("J5").Value2 = 1.8
$ws.Range("L5").Value2 = 12
$ws.Range("N5").Value2 = 8
$ws.Range("Q5").Value2 = 2.1
$ws.Range("R5").Value2 = 1.7
$ws.Range("X5").Value2 = 4.75
$ws.Range("AC5").Value2 = 8
$ws.Range("AD5").Value2 = 10
$ws.Range("AQ5").Value2 = 17
$ws.Range("AW5").Value2 = 12

# Row 8
$ws.Range("G8").Value2 = 1.5
$ws.Range("I8").Value2 = 6.5
$ws.Range("J8").Value2 = 2.05
$ws.Range("Q8").Value2 = 1.83
$ws.Range("R8").Value2 = 2.03
$ws.Range("X8").Value2 = 7
$ws.Range("AH8").Value2 = 17
$ws.Range("AI8").Value2 = 34
$ws.Range("AJ8").Value2 = 21
$ws.Range("AM8").Value2 = 51
$ws.Range("AO8").Value2 = 7.5
$ws.Range("AX8").Value2 = 34
$ws.Range("AZ8").Value2 = 126

# Remove the Danubio vs CA Cerro fixture (id Ei0sKuIl); Excel shifts rows 14-15 up to 13-14
$ws.Rows(13).Delete()

# Two odds corrected on the fixture now at row 13 (Boston River vs Rampla Juniors)
$ws.Range("M13").Value2 = 1.04
$ws.Range("O13").Value2 = 1.3
